$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 (FirstName Hasanul -> Shohana) ---
$ws.Range("B5").Value = "Shohana"

# --- Add two new rows of data ---
$ws.Range("A6").Value = "Ahmed"
$ws.Range("B6").Value = "Munna "
$ws.Range("C6").Value = "Dallas"

$ws.Range("A7").Value = "Touhid"
$ws.Range("B7").Value = "Chowdhury"
$ws.Range("C7").Value = "New Jersey"

# --- Re-add Chowdhury/Sumayla/Calgary and Islam/Shohana/New York rows again ---
$ws.Range("A8").Value = "Chowdhury"
$ws.Range("B8").Value = "Sumayla"
$ws.Range("C8").Value = "Calgary"

$ws.Range("A9").Value = "Islam"
$ws.Range("B9").Value = "Shohana"
$ws.Range("C9").Value = "New York"

# Give C8 the same date-like number format that C4 has (numFmtId 16)
$ws.Range("C4:C8").NumberFormat = "d-mmm"

# --- Borders ---
# Full box (all four sides) around the header + all "plain" data rows
$ws.Range("A1:C5").Borders.LineStyle = 1
$ws.Range("A8:C9").Borders.LineStyle = 1

# Left/Right only borders around the two newly inserted rows
$ws.Range("A6:C7").Borders(7).LineStyle = 1
$ws.Range("A6:C7").Borders(10).LineStyle = 1

# --- Column widths: extend the existing custom width to columns B and C ---
$ws.Range("A1:C1").Columns.AutoFit()
$ws.Columns("B:C").ColumnWidth = $ws.Columns("A").ColumnWidth

# --- Selection matches the author's last active cell ---
$ws.Range("L9").Select()
